$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date-like text to be stored as a literal string (not parsed as a date)
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "09/07/2025"
# Reset the cell style back to the default so it matches the other rows
$ws.Range("A6").Style = "Normal"

$ws.Range("B6").Value = 0.1257222453734942
$ws.Range("C6").Value = 0.8742777546265058
